$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Weekly driver report update: add newly observed driver rows to both the
# "Bad Drivers" and "Good Drivers" tables, and refresh the Totals row.
# ---------------------------------------------------------------------------

# 1) Make room for two new "Bad Driver" rows right above the existing
#    "Totals:" row (old row 4). Excel's Insert() shifts row 4 and everything
#    below it down by 2, carrying formatting from the row above into the
#    freshly inserted rows.
$ws.Rows.Item(4).Resize(2).Insert()

# 2) Make room for three new "Good Driver" rows right below the existing
#    Killer(R) driver row, which is now at row 14 after the shift above.
$ws.Rows.Item(15).Resize(3).Insert()

# ---------------------------------------------------------------------------
# Bad Drivers table
# ---------------------------------------------------------------------------

# Row 3: updated Critical Minutes / Good Roaming Calculation for the existing
# Killer(R) Wi-Fi driver.
$ws.Range("C3").Value = 1247
$ws.Range("D3").Value = 96.5

# Row 4 (new): Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.0.7
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.0.7"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 97.7

# Row 5 (new): Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2"
$ws.Range("B5").Value = 22
$ws.Range("C5").Value = 578
$ws.Range("D5").Value = 98.90000000000001

# Row 6: Totals (was row 4 prior to the insert) - refreshed sums.
$ws.Range("B6").Value = 32
$ws.Range("C6").Value = 1834

# ---------------------------------------------------------------------------
# Good Drivers table
# ---------------------------------------------------------------------------

# Row 15 (new): Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B15").Value = 11140
$ws.Range("D15").Value = 100
# Driver Vintage is stored as plain text (e.g. "2022-08-29"), not a date
# serial, so force a text number format before assigning it - otherwise
# Excel auto-converts the recognizable date pattern to a date value.
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2022-08-29"

# Row 16 (new): Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B16").Value = 14487
$ws.Range("D16").Value = 100
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2022-05-23"

# Row 17 (new): Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B17").Value = 265400
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2022-05-01"
